# Natmi following Dr Hou advice
#
# Refresh the computed NATMI edge-weight statistics for the C3 -> Itgax
# ligand-receptor sheet. Existing sending/target-cluster combinations
# (ECs, FAPs x M2) are recomputed, the sCs target-cluster rows are
# restored, and the new M2/sCs sending-cluster rows are appended so every
# sending cluster x target cluster combination (M2, sCs) is represented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sending cluster, Target cluster labels for rows 2..9, in order.
$sending = @("ECs","ECs","FAPs","FAPs","M2","M2","sCs","sCs")
$target  = @("M2","sCs","M2","sCs","M2","sCs","M2","sCs")

# Numeric columns E..T (20 values per row) for rows 2..9, in order.
$numbers = @(
    # Row 2: ECs -> M2
    @(3,1,1.558564,4.675692,0.005692101168584756,0.005692101168584756,3,1,45.39437433333333,136.183123,0.9965580896473674,0.9965580896473674,70.75003763845733,636.750338746116,0.005672509466644372,0.005672509466644372),
    # Row 3: ECs -> sCs
    @(3,1,1.558564,4.675692,0.005692101168584756,0.005692101168584756,1,0.3333333333333333,0.156783,0.470349,0.00344191035263268,0.00344191035263268,0.244356339612,2.199207056508,[double]"1.959170194038445E-05",[double]"1.959170194038445E-05"),
    # Row 4: FAPs -> M2
    @(3,1,259.5505726666667,778.6517180000001,0.9479162344201305,0.9479162344201304,3,1,45.39437433333333,136.183123,0.9965580896473674,0.9965580896473674,11782.1358540617,106039.2226865553,0.9446535917194513,0.9446535917194512),
    # Row 5: FAPs -> sCs
    @(3,1,259.5505726666667,778.6517180000001,0.9479162344201305,0.9479162344201304,1,0.3333333333333333,0.156783,0.470349,0.00344191035263268,0.00344191035263268,40.69311743439801,366.238056909582,0.003262642700679233,0.003262642700679233),
    # Row 6: M2 -> M2
    @(3,1,12.18925266666667,36.567758,0.04451691386950307,0.04451691386950307,3,1,45.39437433333333,136.183123,0.9965580896473674,0.9965580896473674,553.3234983942482,4979.911485548234,0.04436369064278838,0.04436369064278838),
    # Row 7: M2 -> sCs
    @(3,1,12.18925266666667,36.567758,0.04451691386950307,0.04451691386950307,1,0.3333333333333333,0.156783,0.470349,0.00344191035263268,0.00344191035263268,1.911067600838,17.199608407542,0.0001532232267147,0.0001532232267147),
    # Row 8: sCs -> M2
    @(3,1,0.5133286666666667,1.539986,0.001874750541781658,0.001874750541781658,3,1,45.39437433333333,136.183123,0.9965580896473674,0.9965580896473674,23.30223365069756,209.720102856278,0.001868297818483296,0.001868297818483296),
    # Row 9: sCs -> sCs
    @(3,1,0.5133286666666667,1.539986,0.001874750541781658,0.001874750541781658,1,0.3333333333333333,0.156783,0.470349,0.00344191035263268,0.00344191035263268,0.080481208346,0.724330875114,[double]"6.452723298362014E-06",[double]"6.452723298362014E-06")
)

for ($i = 0; $i -lt 8; $i++) {
    $r = $i + 2

    $ws.Cells.Item($r, 1).Value = $sending[$i]
    $ws.Cells.Item($r, 2).Value = "C3"
    $ws.Cells.Item($r, 3).Value = "Itgax"
    $ws.Cells.Item($r, 4).Value = $target[$i]

    $rowVals = $numbers[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $ws.Cells.Item($r, $j + 5).Value = $rowVals[$j]
    }
}
